# STEP 8 + STEP 9. Final.
# Update the PAS/PAS_CRD passenger figures in column E and refresh the
# current selection (select the whole column F, as last left by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- PAS block (rows 2-9) ---
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 144
$ws.Range("E6").Value = 216
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 504
$ws.Range("E9").Value = 216

# E6 previously carried an explicit "#,##0" number format (style index 1);
# the update drops that custom formatting so the cell reverts to the
# default/general style.
$ws.Range("E6").ClearFormats()

# --- PAS_CRD block (rows 10-17) ---
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 17280
$ws.Range("E14").Value = 25920
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 60480
$ws.Range("E17").Value = 25920

# Leave the sheet with column F selected (matches the saved view state).
$ws.Columns("F").Select()
